# Zoom Recording Links C-DAC.xlsx - add the "OOPJ" and "General Aptitude"
# sections under the existing Logic Building / COS sections.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- xlPasteFormats constant used with PasteSpecial below ---
$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Fix up the two existing section headers (numbering + text tweaks)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "0. Logic Building 20-08-2024 to 24-08-2024"
$ws.Range("A10").Value = "1. Concepts of Operating System (COS) 27-08-2024 to 31-08-2024"

# ---------------------------------------------------------------------
# 2. New section header row 20: "2. Object Oriented Programming in Java (OOPJ)"
# ---------------------------------------------------------------------
$ws.Range("A10:B10").Copy()
$ws.Range("A20:B20").PasteSpecial($xlPasteFormats)
$ws.Rows.Item(20).RowHeight = 15.6
$ws.Range("A20:B20").Merge()
$ws.Range("A20").Value = "2. Object Oriented Programming in Java (OOPJ) 03-09-2024 to 23-09-2024"

# Data rows for the OOPJ section (dates + zoom links)
$oopj = @(
    @{ Row = 21; Date = 45538; Link = "https://us02web.zoom.us/rec/share/Mk0C7rkEvuOAvY-LgYzKHpWu4GaKGJ7gsXbwL12zcKgSlEvdmpFbUdF_N4bHdPpZ.4gopIraI6GVIw5lb" },
    @{ Row = 22; Date = 45539; Link = "https://us02web.zoom.us/rec/share/ctP6sK27-beN-PyhLJ0WAHYVxA1m6Vk56T8ZOsd9keahEVn2tWotCIbF3QcEicGn.kbL89e_IiryJWGhS" },
    @{ Row = 23; Date = 45540; Link = "https://us02web.zoom.us/rec/share/mVQh8dmsxb_wft3eEV5E-lBBRaUjePMF9_f0eRNBLsMdmtVDPpWh_EgJFS3yy-ue.iQkeT6gOm5AHJ6Ys" },
    @{ Row = 24; Date = $null; Link = "https://us02web.zoom.us/rec/share/X7DOTt2ZsZes5PG1CbC3DaXYN0mz4TAgDit7f1_8QyhQLr7WN8HCoNWSY4yUpeo9.StuLQnn24ecEHXbk" },
    @{ Row = 25; Date = 45541; Link = "https://us02web.zoom.us/rec/share/TtKW5GWS0rh4koBuqPCY0iayC6AA09GIp-qQVZkajmRBpHy6j992LvNtdq9yBJY.ZWUzcEg2_gi18d6k" },
    @{ Row = 26; Date = $null; Link = "https://us02web.zoom.us/rec/share/WZc-KhRq40Y8oiUtP1k9N4dvU9Xn3UIwcF5FdxTAqUTp9uODQVfOxTGN-jidrOuJ.2XZYsiXkul8BdYsl" },
    @{ Row = 27; Date = 45542; Link = "https://us02web.zoom.us/rec/share/gSpB5tnOLbY5S44asrpf1yGVX1V8HuRa-IsKEyZfiTE1yukiJtW_efwm3MIY743K.o4haitAfndFseQzq" }
)

foreach ($item in $oopj) {
    $r = $item.Row
    if ($null -ne $item.Date) {
        $ws.Range("A$r").Value = $item.Date
    }
    $ws.Range("B$r").Value = $item.Link
    $ws.Hyperlinks.Add($ws.Range("B$r"), $item.Link)
}

# Re-apply the correct number/hyperlink formatting (the above Hyperlinks.Add
# calls reset the font/style of the target cells to the built-in Hyperlink
# style, so copy the formatting that the rest of the sheet already uses).
$ws.Range("A3").Copy()
$ws.Range("A21").PasteSpecial($xlPasteFormats)
$ws.Range("A22").PasteSpecial($xlPasteFormats)
$ws.Range("A23").PasteSpecial($xlPasteFormats)
$ws.Range("A25").PasteSpecial($xlPasteFormats)
$ws.Range("A27").PasteSpecial($xlPasteFormats)

$ws.Range("B3").Copy()
$ws.Range("B21").PasteSpecial($xlPasteFormats)
$ws.Range("B22").PasteSpecial($xlPasteFormats)
$ws.Range("B23").PasteSpecial($xlPasteFormats)
$ws.Range("B24").PasteSpecial($xlPasteFormats)
$ws.Range("B25").PasteSpecial($xlPasteFormats)
$ws.Range("B26").PasteSpecial($xlPasteFormats)
$ws.Range("B27").PasteSpecial($xlPasteFormats)

# Re-set the dates/links since the PasteSpecial above overwrote cell content
foreach ($item in $oopj) {
    $r = $item.Row
    if ($null -ne $item.Date) {
        $ws.Range("A$r").Value = $item.Date
    }
}

# ---------------------------------------------------------------------
# 3. Rows 28-29: two more links for OOPJ, but NOT wired up as hyperlinks
#    (same visual style as the others, just plain text containing a URL)
# ---------------------------------------------------------------------
$ws.Range("A3").Copy()
$ws.Range("A28").PasteSpecial($xlPasteFormats)
$ws.Range("A28").Value = 45543

$ws.Range("B3").Copy()
$ws.Range("B28").PasteSpecial($xlPasteFormats)
$ws.Range("B28").Value = "https://us02web.zoom.us/rec/share/2Y6CTdRPK0yfPtZ7v3QYVBKEwwwxqiDLofD5BnWNujpMcITWWXxWmt-IbvUZxcQX.PL1Gbn64cZ1VeQ86"

$ws.Range("B3").Copy()
$ws.Range("B29").PasteSpecial($xlPasteFormats)
$ws.Range("B29").Value = "https://us02web.zoom.us/rec/share/A5tgQkWcx2pgU_l_r69DpKtAmAcKyUPg3Rr74DVFsTvKzzHqQ1pR6rFQzG_DVbly.zXM9xSVNsTr5tDio"

# Empty placeholder cell (same style as the link column, but no content)
$ws.Range("B3").Copy()
$ws.Range("B30").PasteSpecial($xlPasteFormats)
$ws.Range("B30").ClearContents()

# ---------------------------------------------------------------------
# 4. New section header row 32: "10. General Aptitude"
# ---------------------------------------------------------------------
$ws.Range("A10:B10").Copy()
$ws.Range("A32:B32").PasteSpecial($xlPasteFormats)
$ws.Rows.Item(32).RowHeight = 15.6
$ws.Range("A32:B32").Merge()
$ws.Range("A32").Value = "10. General Aptitude"

# Row 33: just a date, no link yet
$ws.Range("A3").Copy()
$ws.Range("A33").PasteSpecial($xlPasteFormats)
$ws.Range("A33").Value = 45543

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 5. Column widths
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 40.0534
$ws.Columns.Item(2).ColumnWidth = 70.8307

# ---------------------------------------------------------------------
# 6. View state: zoom + selection
# ---------------------------------------------------------------------
$ws.Range("B30").Select()
$excel.ActiveWindow.Zoom = 70
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
